$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 2 ("Refactoring..."), shifting
# all the existing task rows down by one.
$ws.Rows("2:2").Insert()

# Populate the new row with the "Specular lighting" task (Phong exponent
# export work) and its estimate.
$ws.Range("A2").Value = "Specular lighting"
$ws.Range("B2").Value = 4

# Give the new row its own (non-bold) cell style distinct from the
# default style used by the rest of the data rows.
$ws.Range("A2:B2").Font.Bold = $false
$ws.Range("A2:B2").Font.Name = "Calibri"

# Match the new selection left behind in the saved workbook.
$null = $ws.Range("A2:B2").Select()
